# Regenerate save_data: recompute column G ("K") values from the
# underlying strikeout calc (replacing the old "Strike#" derived values)
# and write the recalculated s_vals back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..68 (one per game row), in order.
$kValues = @(
    0,1,0,0,1,0,0,2,0,0,
    1,1,1,1,2,3,2,1,0,1,
    1,2,1,0,0,0,1,1,1,0,
    3,0,2,1,2,2,1,2,0,1,
    1,1,1,0,1,3,0,0,2,4,
    0,1,2,1,3,3,1,3,1,1,
    0,0,1,2,1,1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
